# Scheduled runner update: refresh market-price-derived columns
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
#  LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ) on the Leve
# profit tables across all job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 950.8
$ws.Range("I19").Value = 755
$ws.Range("K19").Value = 755
$ws.Range("M19").Value = -580
$ws.Range("H28").Value = 598.625
$ws.Range("I28").Value = 555.5714
$ws.Range("K28").Value = 555.5714
$ws.Range("M28").Value = -70.57140000000004
$ws.Range("H43").Value = 3197
$ws.Range("I43").Value = 3197
$ws.Range("K43").Value = 3197
$ws.Range("M43").Value = -3128
$ws.Range("H74").Value = 4000
$ws.Range("I74").Value = 4000
$ws.Range("K74").Value = 4000
$ws.Range("M74").Value = -3064
$ws.Range("H77").Value = 4000
$ws.Range("I77").Value = 4000
$ws.Range("K77").Value = 20000
$ws.Range("M77").Value = -15320
$ws.Range("H86").Value = 3121
$ws.Range("I86").Value = 3026.25
$ws.Range("J86").Value = 3500
$ws.Range("K86").Value = 3026.25
$ws.Range("L86").Value = 3500
$ws.Range("M86").Value = -1903.25
$ws.Range("N86").Value = -5746
$ws.Range("H89").Value = 3121
$ws.Range("I89").Value = 3026.25
$ws.Range("J89").Value = 3500
$ws.Range("K89").Value = 15131.25
$ws.Range("L89").Value = 17500
$ws.Range("M89").Value = -9515.25
$ws.Range("N89").Value = -28732
$ws.Range("H92").Value = 128.5
$ws.Range("I92").Value = 128.5
$ws.Range("K92").Value = 128.5
$ws.Range("M92").Value = 1119.5
$ws.Range("H137").Value = 2499.25
$ws.Range("I137").Value = 1999.6666
$ws.Range("J137").Value = 3998
$ws.Range("K137").Value = 5998.9998
$ws.Range("L137").Value = 11994
$ws.Range("M137").Value = -3448.9998
$ws.Range("N137").Value = -17094
$ws.Range("H138").Value = 8227.16
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 8227.16
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 24681.48
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = -34961.48

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3455
$ws.Range("I32").Value = 3581.5652
$ws.Range("K32").Value = 3581.5652
$ws.Range("M32").Value = -3294.5652
$ws.Range("H61").Value = 7900.6
$ws.Range("I61").Value = 9125.75
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 9125.75
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -8913.75
$ws.Range("N61").Value = -3424
$ws.Range("H136").Value = 7900.6
$ws.Range("I136").Value = 9125.75
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 27377.25
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -24827.25
$ws.Range("N136").Value = -14100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3381.1333
$ws.Range("I20").Value = 3327.9167
$ws.Range("K20").Value = 3327.9167
$ws.Range("M20").Value = -3080.9167
$ws.Range("H105").Value = 1535.2727
$ws.Range("I105").Value = 1535.2727
$ws.Range("K105").Value = 1535.2727
$ws.Range("M105").Value = 211.7273
$ws.Range("H134").Value = 3639.8
$ws.Range("I134").Value = 3578.3572
$ws.Range("K134").Value = 10735.0716
$ws.Range("M134").Value = -8200.071599999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 8658.75
$ws.Range("J16").Value = 8531
$ws.Range("L16").Value = 8531
$ws.Range("N16").Value = -9105
$ws.Range("H31").Value = 10598.8
$ws.Range("I31").Value = 1666.3334
$ws.Range("J31").Value = 23997.5
$ws.Range("K31").Value = 1666.3334
$ws.Range("L31").Value = 23997.5
$ws.Range("M31").Value = -1371.3334
$ws.Range("N31").Value = -24587.5
$ws.Range("H34").Value = 10598.8
$ws.Range("I34").Value = 1666.3334
$ws.Range("J34").Value = 23997.5
$ws.Range("K34").Value = 1666.3334
$ws.Range("L34").Value = 23997.5
$ws.Range("M34").Value = -1464.3334
$ws.Range("N34").Value = -24401.5
$ws.Range("H113").Value = 8658.75
$ws.Range("J113").Value = 8531
$ws.Range("L113").Value = 8531
$ws.Range("N113").Value = -12871

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 11028.143
$ws.Range("I80").Value = 10800
$ws.Range("J80").Value = 11066.167
$ws.Range("K80").Value = 32400
$ws.Range("L80").Value = 33198.501
$ws.Range("N80").Value = -35070.501
$ws.Range("M80").Value = -31464
$ws.Range("H83").Value = 11028.143
$ws.Range("I83").Value = 10800
$ws.Range("J83").Value = 11066.167
$ws.Range("K83").Value = 97200
$ws.Range("L83").Value = 99595.503
$ws.Range("N83").Value = -108955.503
$ws.Range("M83").Value = -92520
$ws.Range("H104").Value = 4685.7144
$ws.Range("J104").Value = 4685.7144
$ws.Range("L104").Value = 14057.1432
$ws.Range("N104").Value = -19299.1432
$ws.Range("H131").Value = 1533.1666
$ws.Range("I131").Value = 1133
$ws.Range("J131").Value = 1933.3334
$ws.Range("K131").Value = 3399
$ws.Range("L131").Value = 5800.0002
$ws.Range("M131").Value = 1641
$ws.Range("N131").Value = -15880.0002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3087.889
$ws.Range("I102").Value = 3087.889
$ws.Range("K102").Value = 3087.889
$ws.Range("M102").Value = -1465.889

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 5731.3335
$ws.Range("I46").Value = 9999
$ws.Range("J46").Value = 3597.5
$ws.Range("K46").Value = 9999
$ws.Range("L46").Value = 3597.5
$ws.Range("N46").Value = -3973.5
$ws.Range("M46").Value = -9811
$ws.Range("H55").Value = 50.666668
$ws.Range("I55").Value = 56
$ws.Range("J55").Value = 40
$ws.Range("K55").Value = 56
$ws.Range("L55").Value = 40
$ws.Range("M55").Value = 117
$ws.Range("N55").Value = -386
$ws.Range("H61").Value = 3540.8333
$ws.Range("I61").Value = 3499
$ws.Range("K61").Value = 3499
$ws.Range("M61").Value = -3297
$ws.Range("H113").Value = 3540.8333
$ws.Range("I113").Value = 3499
$ws.Range("K113").Value = 3499
$ws.Range("M113").Value = -1329

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5917
$ws.Range("I62").Value = 5917
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 5917
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -5293
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 5917
$ws.Range("I65").Value = 5917
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 29585
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -26465
$ws.Range("N65").ClearContents()
